# Update cryptos list with new price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay text even when the new value looks like a plain number
    # (e.g. "21.69"), then restore the cell's style back to the default (Normal) so
    # no stray number-format styling is left behind on the cell.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.873.86"
$ws.Range("E2").Value = "  +0.82%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.551.82"
$ws.Range("E3").Value = "  +1.52%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.39%  "

# Row 5 - BNB
Set-TextValue "D5" "206.87"
$ws.Range("E5").Value = "  +0.81%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.34%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.35%  "

# Row 8 - Solana
Set-TextValue "D8" "21.69"
$ws.Range("E8").Value = "  +2.09%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.50%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.43%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.89%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.774.71"
$ws.Range("E12").Value = "  +1.54%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.559.94"
$ws.Range("E13").Value = "  +1.93%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.514"
$ws.Range("E15").Value = "  +2.04%  "

# Row 16 & 17 swap: Litecoin <-> WrappedBTC
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "26.879.59"
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D17" "61.69"
$ws.Range("E17").Value = "  +0.78%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "216.58"
$ws.Range("E18").Value = "  +2.43%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +0.87%  "

# Row 20 - Chainlink
Set-TextValue "D20" "7.21"
$ws.Range("E20").Value = "  +0.82%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.12%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.17"
$ws.Range("E23").Value = "  +1.52%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +2.39%  "

# Row 25 - Monero
Set-TextValue "D25" "153.55"
$ws.Range("E25").Value = "  +0.90%  "

# Row 26 - Cosmos
Set-TextValue "D26" "6.63"
$ws.Range("E26").Value = "  +1.80%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "14.86"
$ws.Range("E27").Value = "  +0.36%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.38%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  +1.50%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.09"
$ws.Range("E31").Value = "  -0.07%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.08%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.417.74"
$ws.Range("E33").Value = "  +4.66%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +3.41%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "1.56"
$ws.Range("E35").Value = "  +4.43%  "

# Row 36 - TrustWalletToken
Set-TextValue "D36" "0.959"
$ws.Range("E36").Value = "  +2.44%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.71%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +1.22%  "

# Row 39 - ImmutableX
Set-TextValue "D39" "0.523"
$ws.Range("E39").Value = "  +0.30%  "

# Row 40 - ARBITRUM
Set-TextValue "D40" "0.806"
$ws.Range("E40").Value = "  +1.50%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.39%  "

# Row 42 - FraxShare
Set-TextValue "D42" "5.63"
$ws.Range("E42").Value = "  -1.06%  "

# Row 43 - WEMIXToken
$ws.Range("E43").Value = "  -0.64%  "

# Row 44 - MXToken
Set-TextValue "D44" "2.26"
$ws.Range("E44").Value = "  +3.91%  "

# Row 45 - Aave
Set-TextValue "D45" "63.62"
$ws.Range("E45").Value = "  +2.20%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +1.44%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.688.31"
$ws.Range("E47").Value = "  +1.44%  "

# Row 48 - Quant
Set-TextValue "D48" "86.05"
$ws.Range("E48").Value = "  +0.50%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0518"
$ws.Range("E49").Value = "  +1.99%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +5.89%  "

# Row 51 - Algorand
Set-TextValue "D51" "0.0962"
$ws.Range("E51").Value = "  +1.93%  "
